$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Ccl12"
$ws.Cells.Item(2,3).Value = "Ccr5"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = [double]"1"
$ws.Cells.Item(2,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(2,7).Value = [double]"0.2821546666666667"
$ws.Cells.Item(2,8).Value = [double]"0.846464"
$ws.Cells.Item(2,9).Value = [double]"0.003251136722938652"
$ws.Cells.Item(2,10).Value = [double]"0.003251136722938651"
$ws.Cells.Item(2,11).Value = [double]"3"
$ws.Cells.Item(2,12).Value = [double]"1"
$ws.Cells.Item(2,13).Value = [double]"0.1207436666666667"
$ws.Cells.Item(2,14).Value = [double]"0.362231"
$ws.Cells.Item(2,15).Value = [double]"0.001088347656299082"
$ws.Cells.Item(2,16).Value = [double]"0.001088347656299082"
$ws.Cells.Item(2,17).Value = [double]"0.03406838902044444"
$ws.Cells.Item(2,18).Value = [double]"0.306615501184"
$ws.Cells.Item(2,19).Value = [double]"3.53836703271816E-06"
$ws.Cells.Item(2,20).Value = [double]"3.53836703271816E-06"

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Ccl12"
$ws.Cells.Item(3,3).Value = "Ccr5"
$ws.Cells.Item(3,4).Value = "M1"
$ws.Cells.Item(3,5).Value = [double]"1"
$ws.Cells.Item(3,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(3,7).Value = [double]"0.2821546666666667"
$ws.Cells.Item(3,8).Value = [double]"0.846464"
$ws.Cells.Item(3,9).Value = [double]"0.003251136722938652"
$ws.Cells.Item(3,10).Value = [double]"0.003251136722938651"
$ws.Cells.Item(3,11).Value = [double]"3"
$ws.Cells.Item(3,12).Value = [double]"1"
$ws.Cells.Item(3,13).Value = [double]"48.86031499999999"
$ws.Cells.Item(3,14).Value = [double]"146.580945"
$ws.Cells.Item(3,15).Value = [double]"0.4404124107236948"
$ws.Cells.Item(3,16).Value = [double]"0.4404124107236948"
$ws.Cells.Item(3,17).Value = [double]"13.78616589205333"
$ws.Cells.Item(3,18).Value = [double]"124.07549302848"
$ws.Cells.Item(3,19).Value = [double]"0.001431840961741745"
$ws.Cells.Item(3,20).Value = [double]"0.001431840961741744"

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Ccl12"
$ws.Cells.Item(4,3).Value = "Ccr5"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = [double]"1"
$ws.Cells.Item(4,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(4,7).Value = [double]"0.2821546666666667"
$ws.Cells.Item(4,8).Value = [double]"0.846464"
$ws.Cells.Item(4,9).Value = [double]"0.003251136722938652"
$ws.Cells.Item(4,10).Value = [double]"0.003251136722938651"
$ws.Cells.Item(4,11).Value = [double]"3"
$ws.Cells.Item(4,12).Value = [double]"1"
$ws.Cells.Item(4,13).Value = [double]"61.90825266666667"
$ws.Cells.Item(4,14).Value = [double]"185.724758"
$ws.Cells.Item(4,15).Value = [double]"0.5580226570503747"
$ws.Cells.Item(4,16).Value = [double]"0.5580226570503747"
$ws.Cells.Item(4,17).Value = [double]"17.46770239507911"
$ws.Cells.Item(4,18).Value = [double]"157.209321555712"
$ws.Cells.Item(4,19).Value = [double]"0.001814207952568274"
$ws.Cells.Item(4,20).Value = [double]"0.001814207952568274"

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ccl12"
$ws.Cells.Item(5,3).Value = "Ccr5"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = [double]"1"
$ws.Cells.Item(5,6).Value = [double]"0.3333333333333333"
$ws.Cells.Item(5,7).Value = [double]"0.2821546666666667"
$ws.Cells.Item(5,8).Value = [double]"0.846464"
$ws.Cells.Item(5,9).Value = [double]"0.003251136722938652"
$ws.Cells.Item(5,10).Value = [double]"0.003251136722938651"
$ws.Cells.Item(5,11).Value = [double]"2"
$ws.Cells.Item(5,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(5,13).Value = [double]"0.05287333333333333"
$ws.Cells.Item(5,14).Value = [double]"0.15862"
$ws.Cells.Item(5,15).Value = [double]"0.0004765845696314243"
$ws.Cells.Item(5,16).Value = [double]"0.0004765845696314243"
$ws.Cells.Item(5,17).Value = [double]"0.01491845774222222"
$ws.Cells.Item(5,18).Value = [double]"0.13426611968"
$ws.Cells.Item(5,19).Value = [double]"1.549441595914636E-06"
$ws.Cells.Item(5,20).Value = [double]"1.549441595914636E-06"

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Ccl12"
$ws.Cells.Item(6,3).Value = "Ccr5"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = [double]"3"
$ws.Cells.Item(6,6).Value = [double]"1"
$ws.Cells.Item(6,7).Value = [double]"55.96940366666666"
$ws.Cells.Item(6,8).Value = [double]"167.908211"
$ws.Cells.Item(6,9).Value = [double]"0.6449093533393406"
$ws.Cells.Item(6,10).Value = [double]"0.6449093533393405"
$ws.Cells.Item(6,11).Value = [double]"3"
$ws.Cells.Item(6,12).Value = [double]"1"
$ws.Cells.Item(6,13).Value = [double]"0.1207436666666667"
$ws.Cells.Item(6,14).Value = [double]"0.362231"
$ws.Cells.Item(6,15).Value = [double]"0.001088347656299082"
$ws.Cells.Item(6,16).Value = [double]"0.001088347656299082"
$ws.Cells.Item(6,17).Value = [double]"6.75795101986011"
$ws.Cells.Item(6,18).Value = [double]"60.82155917874099"
$ws.Cells.Item(6,19).Value = [double]"0.0007018855832322281"
$ws.Cells.Item(6,20).Value = [double]"0.000701885583232228"

# Row 7
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Ccl12"
$ws.Cells.Item(7,3).Value = "Ccr5"
$ws.Cells.Item(7,4).Value = "M1"
$ws.Cells.Item(7,5).Value = [double]"3"
$ws.Cells.Item(7,6).Value = [double]"1"
$ws.Cells.Item(7,7).Value = [double]"55.96940366666666"
$ws.Cells.Item(7,8).Value = [double]"167.908211"
$ws.Cells.Item(7,9).Value = [double]"0.6449093533393406"
$ws.Cells.Item(7,10).Value = [double]"0.6449093533393405"
$ws.Cells.Item(7,11).Value = [double]"3"
$ws.Cells.Item(7,12).Value = [double]"1"
$ws.Cells.Item(7,13).Value = [double]"48.86031499999999"
$ws.Cells.Item(7,14).Value = [double]"146.580945"
$ws.Cells.Item(7,15).Value = [double]"0.4404124107236948"
$ws.Cells.Item(7,16).Value = [double]"0.4404124107236948"
$ws.Cells.Item(7,17).Value = [double]"2734.682693515488"
$ws.Cells.Item(7,18).Value = [double]"24612.14424163939"
$ws.Cells.Item(7,19).Value = [double]"0.2840260830024381"
$ws.Cells.Item(7,20).Value = [double]"0.284026083002438"

# Row 8
$ws.Cells.Item(8,1).Value = "M1"
$ws.Cells.Item(8,2).Value = "Ccl12"
$ws.Cells.Item(8,3).Value = "Ccr5"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = [double]"3"
$ws.Cells.Item(8,6).Value = [double]"1"
$ws.Cells.Item(8,7).Value = [double]"55.96940366666666"
$ws.Cells.Item(8,8).Value = [double]"167.908211"
$ws.Cells.Item(8,9).Value = [double]"0.6449093533393406"
$ws.Cells.Item(8,10).Value = [double]"0.6449093533393405"
$ws.Cells.Item(8,11).Value = [double]"3"
$ws.Cells.Item(8,12).Value = [double]"1"
$ws.Cells.Item(8,13).Value = [double]"61.90825266666667"
$ws.Cells.Item(8,14).Value = [double]"185.724758"
$ws.Cells.Item(8,15).Value = [double]"0.5580226570503747"
$ws.Cells.Item(8,16).Value = [double]"0.5580226570503747"
$ws.Cells.Item(8,17).Value = [double]"3464.96798379866"
$ws.Cells.Item(8,18).Value = [double]"31184.71185418794"
$ws.Cells.Item(8,19).Value = [double]"0.3598740309070578"
$ws.Cells.Item(8,20).Value = [double]"0.3598740309070577"

# Row 9
$ws.Cells.Item(9,1).Value = "M1"
$ws.Cells.Item(9,2).Value = "Ccl12"
$ws.Cells.Item(9,3).Value = "Ccr5"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = [double]"3"
$ws.Cells.Item(9,6).Value = [double]"1"
$ws.Cells.Item(9,7).Value = [double]"55.96940366666666"
$ws.Cells.Item(9,8).Value = [double]"167.908211"
$ws.Cells.Item(9,9).Value = [double]"0.6449093533393406"
$ws.Cells.Item(9,10).Value = [double]"0.6449093533393405"
$ws.Cells.Item(9,11).Value = [double]"2"
$ws.Cells.Item(9,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(9,13).Value = [double]"0.05287333333333333"
$ws.Cells.Item(9,14).Value = [double]"0.15862"
$ws.Cells.Item(9,15).Value = [double]"0.0004765845696314243"
$ws.Cells.Item(9,16).Value = [double]"0.0004765845696314243"
$ws.Cells.Item(9,17).Value = [double]"2.959288936535556"
$ws.Cells.Item(9,18).Value = [double]"26.63360042882"
$ws.Cells.Item(9,19).Value = [double]"0.0003073538466125098"
$ws.Cells.Item(9,20).Value = [double]"0.0003073538466125098"

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Ccl12"
$ws.Cells.Item(10,3).Value = "Ccr5"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = [double]"3"
$ws.Cells.Item(10,6).Value = [double]"1"
$ws.Cells.Item(10,7).Value = [double]"30.53490766666667"
$ws.Cells.Item(10,8).Value = [double]"91.60472300000001"
$ws.Cells.Item(10,9).Value = [double]"0.3518395099377208"
$ws.Cells.Item(10,10).Value = [double]"0.3518395099377208"
$ws.Cells.Item(10,11).Value = [double]"3"
$ws.Cells.Item(10,12).Value = [double]"1"
$ws.Cells.Item(10,13).Value = [double]"0.1207436666666667"
$ws.Cells.Item(10,14).Value = [double]"0.362231"
$ws.Cells.Item(10,15).Value = [double]"0.001088347656299082"
$ws.Cells.Item(10,16).Value = [double]"0.001088347656299082"
$ws.Cells.Item(10,17).Value = [double]"3.686896713001444"
$ws.Cells.Item(10,18).Value = [double]"33.182070417013"
$ws.Cells.Item(10,19).Value = [double]"0.0003829237060341361"
$ws.Cells.Item(10,20).Value = [double]"0.0003829237060341361"

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Ccl12"
$ws.Cells.Item(11,3).Value = "Ccr5"
$ws.Cells.Item(11,4).Value = "M1"
$ws.Cells.Item(11,5).Value = [double]"3"
$ws.Cells.Item(11,6).Value = [double]"1"
$ws.Cells.Item(11,7).Value = [double]"30.53490766666667"
$ws.Cells.Item(11,8).Value = [double]"91.60472300000001"
$ws.Cells.Item(11,9).Value = [double]"0.3518395099377208"
$ws.Cells.Item(11,10).Value = [double]"0.3518395099377208"
$ws.Cells.Item(11,11).Value = [double]"3"
$ws.Cells.Item(11,12).Value = [double]"1"
$ws.Cells.Item(11,13).Value = [double]"48.86031499999999"
$ws.Cells.Item(11,14).Value = [double]"146.580945"
$ws.Cells.Item(11,15).Value = [double]"0.4404124107236948"
$ws.Cells.Item(11,16).Value = [double]"0.4404124107236948"
$ws.Cells.Item(11,17).Value = [double]"1491.945207089248"
$ws.Cells.Item(11,18).Value = [double]"13427.50686380324"
$ws.Cells.Item(11,19).Value = [double]"0.154954486759515"
$ws.Cells.Item(11,20).Value = [double]"0.154954486759515"

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Ccl12"
$ws.Cells.Item(12,3).Value = "Ccr5"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = [double]"3"
$ws.Cells.Item(12,6).Value = [double]"1"
$ws.Cells.Item(12,7).Value = [double]"30.53490766666667"
$ws.Cells.Item(12,8).Value = [double]"91.60472300000001"
$ws.Cells.Item(12,9).Value = [double]"0.3518395099377208"
$ws.Cells.Item(12,10).Value = [double]"0.3518395099377208"
$ws.Cells.Item(12,11).Value = [double]"3"
$ws.Cells.Item(12,12).Value = [double]"1"
$ws.Cells.Item(12,13).Value = [double]"61.90825266666667"
$ws.Cells.Item(12,14).Value = [double]"185.724758"
$ws.Cells.Item(12,15).Value = [double]"0.5580226570503747"
$ws.Cells.Item(12,16).Value = [double]"0.5580226570503747"
$ws.Cells.Item(12,17).Value = [double]"1890.362778981337"
$ws.Cells.Item(12,18).Value = [double]"17013.26501083204"
$ws.Cells.Item(12,19).Value = [double]"0.1963344181907487"
$ws.Cells.Item(12,20).Value = [double]"0.1963344181907486"

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Ccl12"
$ws.Cells.Item(13,3).Value = "Ccr5"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = [double]"3"
$ws.Cells.Item(13,6).Value = [double]"1"
$ws.Cells.Item(13,7).Value = [double]"30.53490766666667"
$ws.Cells.Item(13,8).Value = [double]"91.60472300000001"
$ws.Cells.Item(13,9).Value = [double]"0.3518395099377208"
$ws.Cells.Item(13,10).Value = [double]"0.3518395099377208"
$ws.Cells.Item(13,11).Value = [double]"2"
$ws.Cells.Item(13,12).Value = [double]"0.6666666666666666"
$ws.Cells.Item(13,13).Value = [double]"0.05287333333333333"
$ws.Cells.Item(13,14).Value = [double]"0.15862"
$ws.Cells.Item(13,15).Value = [double]"0.0004765845696314243"
$ws.Cells.Item(13,16).Value = [double]"0.0004765845696314243"
$ws.Cells.Item(13,17).Value = [double]"1.614482351362222"
$ws.Cells.Item(13,18).Value = [double]"14.53034116226"
$ws.Cells.Item(13,19).Value = [double]"0.0001676812814229999"
$ws.Cells.Item(13,20).Value = [double]"0.0001676812814229999"
